$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.854.93'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '1.755.27'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'327.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = "'0.4595"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = "'41.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D14").Value = "'5.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").Value = "'7.157"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '1.754.64'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = "'91.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = "'0.06416"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").Value = "'16.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("D22").Value = "'5.743"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '27.885.60'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("D24").Value = "'11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").Value = "'2.159"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.17%  '
$ws.Range("D26").Value = "'161.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("D27").Value = "'20.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("D28").Value = '1.959.67'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = "'2.144"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = "'123.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("D31").Value = "'1.065"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").Value = "'0.09261"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").Value = "'5.526"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = "'0.02266"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").Value = "'0.06081"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("D38").Value = "'0.2056"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("D39").Value = "'4.889"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").Value = "'1.357"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").Value = "'7.752"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").Value = "'13.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").Value = "'3.723"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = "'0.5777"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").Value = "'123.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").Value = "'1.922"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("D49").Value = "'0.06795"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").Value = "'1.119"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = "'71.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
